# Apply cell value updates per the commit diff (cryptos list refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '64.735.42'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  +5.31%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '3.102.48'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  +3.46%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.00'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  -0.08%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '559.62'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +2.97%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '143.87'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +10.42%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.999'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -0.08%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '3.101.26'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +3.50%  '
$ws.Range('E9').Value = '  +2.11%  '
$ws.Range('E10').Value = '  +17.71%  '
$ws.Range('E11').Value = '  +5.41%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.462'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +4.42%  '
$ws.Range('E13').Value = '  +4.66%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '35.32'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +4.53%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '3.605.42'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +3.40%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '64.701.89'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +4.95%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '3.100.90'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +3.28%  '
$ws.Range('E18').Value = '  -0.45%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '6.84'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +3.96%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '484.55'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +1.15%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '13.86'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +5.52%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '7.65'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +10.31%  '
$ws.Range('E23').Value = '  +2.10%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '13.30'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +11.19%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '80.98'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +0.29%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '1.00'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +0.02%  '
$ws.Range('E27').Value = '  +4.16%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '8.12'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +6.91%  '
$ws.Range('E29').Value = '  +9.15%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.999'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -0.13%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '26.13'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +2.58%  '
$ws.Range('E32').Value = '  +3.37%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '2.46'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +6.12%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '5.72'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +3.77%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '6.21'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +6.57%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '55.07'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +0.78%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '465.17'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +4.92%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.0409'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +7.46%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.0829'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +5.24%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '3.013.33'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -3.55%  '
$ws.Range('E41').Value = '  +1.01%  '
$ws.Range('B42').Value = 'dogwifhat'
$ws.Range('C42').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '2.74'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +16.62%  '
$ws.Range('B43').Value = 'Cosmos'
$ws.Range('C43').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '8.27'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +2.95%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '28.56'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +10.99%  '
$ws.Range('E45').Value = '  +8.07%  '
$ws.Range('E47').Value = '  +9.08%  '
$ws.Range('E48').Value = '  +4.58%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '118.90'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +4.49%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.0₃0516'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +7.34%  '
$ws.Range('E51').Value = '  +3.56%  '
